# Apply updated market-price / profit figures to each job sheet.
# Values sourced from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1599
$ws.Range("I43").Value = 1999
$ws.Range("J43").Value = 1199
$ws.Range("K43").Value = 1999
$ws.Range("L43").Value = 1199
$ws.Range("M43").Value = -1930
$ws.Range("N43").Value = -1337
$ws.Range("H100").Value = 2686.8667
$ws.Range("I100").Value = 2413.182
$ws.Range("J100").Value = 3439.5
$ws.Range("K100").Value = 2413.182
$ws.Range("L100").Value = 3439.5
$ws.Range("M100").Value = -1872.182
$ws.Range("N100").Value = -4521.5
$ws.Range("H101").Value = 763.8889
$ws.Range("I101").Value = 484.375
$ws.Range("K101").Value = 1453.125
$ws.Range("M101").Value = 168.875
$ws.Range("H108").Value = 96623.336
$ws.Range("J108").Value = 96623.336
$ws.Range("L108").Value = 96623.336
$ws.Range("N108").Value = -104303.336
$ws.Range("H109").Value = 49703.57
$ws.Range("J109").Value = 49703.57
$ws.Range("L109").Value = 49703.57
$ws.Range("N109").Value = -52477.57
$ws.Range("H110").Value = 52782.43
$ws.Range("J110").Value = 52782.43
$ws.Range("L110").Value = 52782.43
$ws.Range("N110").Value = -60962.43
$ws.Range("H132").Value = 1601.3334
$ws.Range("I132").Value = 1622.069
$ws.Range("K132").Value = 4866.207
$ws.Range("M132").Value = -2336.207
$ws.Range("H133").Value = 78526.61
$ws.Range("J133").Value = 78526.61
$ws.Range("L133").Value = 78526.61
$ws.Range("N133").Value = -88646.61
$ws.Range("H134").Value = 94216.664
$ws.Range("J134").Value = 94216.664
$ws.Range("L134").Value = 94216.664
$ws.Range("N134").Value = -104356.664
$ws.Range("H136").Value = 99995
$ws.Range("J136").Value = 99995
$ws.Range("L136").Value = 99995
$ws.Range("N136").Value = -110195

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 98990
$ws.Range("J7").Value = 98990
$ws.Range("L7").Value = 98990
$ws.Range("N7").Value = -99218
$ws.Range("H45").Value = 10422168
$ws.Range("I45").Value = 2502.5
$ws.Range("K45").Value = 2502.5
$ws.Range("M45").Value = -2125.5
$ws.Range("H102").Value = 123635.336
$ws.Range("I102").Value = 144370.86
$ws.Range("K102").Value = 144370.86
$ws.Range("M102").Value = -142748.86
$ws.Range("H107").Value = 41205.715
$ws.Range("J107").Value = 41205.715
$ws.Range("L107").Value = 41205.715
$ws.Range("N107").Value = -48885.715
$ws.Range("H127").Value = 84994.336
$ws.Range("J127").Value = 84994.336
$ws.Range("L127").Value = 84994.336
$ws.Range("N127").Value = -94914.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 86988
$ws.Range("J13").Value = 86988
$ws.Range("L13").Value = 86988
$ws.Range("N13").Value = -87324
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H53").Value = 68893
$ws.Range("J53").Value = 68893
$ws.Range("L53").Value = 68893
$ws.Range("N53").Value = -70041
$ws.Range("H94").Value = 3476.8572
$ws.Range("I94").Value = 2865.6
$ws.Range("J94").Value = 5005
$ws.Range("K94").Value = 2865.6
$ws.Range("L94").Value = 5005
$ws.Range("M94").Value = -2414.6
$ws.Range("N94").Value = -5907
$ws.Range("H107").Value = 3177.5
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120
$ws.Range("H108").Value = 99994
$ws.Range("J108").Value = 99994
$ws.Range("L108").Value = 99994
$ws.Range("N108").Value = -107674
$ws.Range("H109").Value = 89659.336
$ws.Range("J109").Value = 89659.336
$ws.Range("L109").Value = 89659.336
$ws.Range("N109").Value = -92433.336
$ws.Range("H114").Value = 75996.664
$ws.Range("J114").Value = 75996.664
$ws.Range("L114").Value = 75996.664
$ws.Range("N114").Value = -84674.664
$ws.Range("H115").Value = 76568.71000000001
$ws.Range("J115").Value = 78996.664
$ws.Range("L115").Value = 78996.664
$ws.Range("N115").Value = -82130.664
$ws.Range("H118").Value = 65163.5
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H122").Value = 77775.336
$ws.Range("J122").Value = 77775.336
$ws.Range("L122").Value = 77775.336
$ws.Range("N122").Value = -87575.336
$ws.Range("H127").Value = 57124.57
$ws.Range("J127").Value = 57124.57
$ws.Range("L127").Value = 57124.57
$ws.Range("N127").Value = -67044.57000000001
$ws.Range("H132").Value = 98321.664
$ws.Range("J132").Value = 98321.664
$ws.Range("L132").Value = 98321.664
$ws.Range("N132").Value = -108441.664
$ws.Range("H138").Value = 89996.664
$ws.Range("J138").Value = 89996.664
$ws.Range("L138").Value = 89996.664
$ws.Range("N138").Value = -100276.664
$ws.Range("H140").Value = 76520
$ws.Range("J140").Value = 89990
$ws.Range("L140").Value = 89990
$ws.Range("N140").Value = -100350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 43128.168
$ws.Range("J9").Value = 43128.168
$ws.Range("L9").Value = 43128.168
$ws.Range("N9").Value = -43464.168
$ws.Range("H58").Value = 1435.8223
$ws.Range("I58").Value = 1315.9032
$ws.Range("J58").Value = 1701.3572
$ws.Range("K58").Value = 1315.9032
$ws.Range("L58").Value = 1701.3572
$ws.Range("M58").Value = -1112.9032
$ws.Range("N58").Value = -2107.3572
$ws.Range("H114").Value = 39611.5
$ws.Range("J114").Value = 39611.5
$ws.Range("L114").Value = 39611.5
$ws.Range("N114").Value = -48289.5
$ws.Range("H118").Value = 59854.285
$ws.Range("J118").Value = 59854.285
$ws.Range("L118").Value = 59854.285
$ws.Range("N118").Value = -63168.285
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H134").Value = 64646.875
$ws.Range("I134").Value = 1578.25
$ws.Range("K134").Value = 4734.75
$ws.Range("M134").Value = -2199.75
$ws.Range("H136").Value = 1435.8223
$ws.Range("I136").Value = 1315.9032
$ws.Range("J136").Value = 1701.3572
$ws.Range("K136").Value = 3947.7096
$ws.Range("L136").Value = 5104.071599999999
$ws.Range("M136").Value = -1397.7096
$ws.Range("N136").Value = -10204.0716
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 7495
$ws.Range("J48").Value = 7495
$ws.Range("L48").Value = 22485
$ws.Range("N48").Value = -22985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20665.445
$ws.Range("J93").Value = 20665.445
$ws.Range("L93").Value = 20665.445
$ws.Range("N93").Value = -24409.445
$ws.Range("H110").Value = 73802.25
$ws.Range("J110").Value = 73802.25
$ws.Range("L110").Value = 73802.25
$ws.Range("N110").Value = -81982.25
$ws.Range("H116").Value = 59854.285
$ws.Range("J116").Value = 59854.285
$ws.Range("L116").Value = 59854.285
$ws.Range("N116").Value = -69032.285
$ws.Range("H119").Value = 57067.5
$ws.Range("J119").Value = 56442.668
$ws.Range("L119").Value = 56442.668
$ws.Range("N119").Value = -66118.66800000001
$ws.Range("H135").Value = 59902.617
$ws.Range("J135").Value = 59902.617
$ws.Range("L135").Value = 59902.617
$ws.Range("N135").Value = -70042.617
$ws.Range("H140").Value = 98567.28999999999
$ws.Range("J140").Value = 98567.28999999999
$ws.Range("L140").Value = 98567.28999999999
$ws.Range("N140").Value = -108927.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1893.3334
$ws.Range("I93").Value = 1893.3334
$ws.Range("K93").Value = 1893.3334
$ws.Range("M93").Value = -645.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 41276
$ws.Range("J121").Value = 41276
$ws.Range("L121").Value = 41276
$ws.Range("N121").Value = -44770
$ws.Range("H132").Value = 13330.519
$ws.Range("I132").Value = 15609.954
$ws.Range("K132").Value = 46829.862
$ws.Range("M132").Value = -44299.862
